$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Price" (column D) and "Volume(1h)" (column E) values ---
#
# The Price column holds plain text in the source workbook (numbers are
# formatted with "." as a thousands separator, e.g. "59.481.59", and some
# values like "1.00" must keep their trailing zeros). Assigning such a
# numeric-looking string straight to a cell's .Value makes Excel silently
# reinterpret it as a real number (losing trailing zeros / exact text).
# To avoid that, the cells that need it are temporarily switched to the
# Text number format ("@") before the value is written, then restored to
# the default "Normal" style so the cell formatting matches the original
# file once the text has been stored.

$priceCells = @(
    'D2',
    'D3',
    'D5',
    'D6',
    'D10',
    'D11',
    'D13',
    'D14',
    'D15',
    'D17',
    'D18',
    'D20',
    'D22',
    'D23',
    'D26',
    'D29',
    'D30',
    'D31',
    'D32',
    'D33',
    'D35',
    'D38',
    'D39',
    'D40',
    'D41',
    'D42',
    'D43',
    'D44',
    'D46',
    'D47',
    'D48',
    'D50',
    'D51'
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range('D2').Value = '59.405.94'
$ws.Range('E2').Value = '  +2.59%  '

# Row 3
$ws.Range('D3').Value = '2.399.66'
$ws.Range('E3').Value = '  +2.69%  '

# Row 5
$ws.Range('D5').Value = '550.57'
$ws.Range('E5').Value = '  +2.14%  '

# Row 6
$ws.Range('D6').Value = '136.28'
$ws.Range('E6').Value = '  +1.50%  '

# Row 7
$ws.Range('E7').Value = '  +0.05%  '

# Row 8
$ws.Range('E8').Value = '  +2.43%  '

# Row 9
$ws.Range('E9').Value = '  +5.72%  '

# Row 10
$ws.Range('D10').Value = '5.79'
$ws.Range('E10').Value = '  +4.86%  '

# Row 11
$ws.Range('D11').Value = '0.360'
$ws.Range('E11').Value = '  +1.19%  '

# Row 12
$ws.Range('E12').Value = '  -2.42%  '

# Row 13
$ws.Range('D13').Value = '24.59'
$ws.Range('E13').Value = '  +3.86%  '

# Row 14
$ws.Range('D14').Value = '2.833.32'
$ws.Range('E14').Value = '  +3.19%  '

# Row 15
$ws.Range('D15').Value = '59.351.95'
$ws.Range('E15').Value = '  +2.60%  '

# Row 16
$ws.Range('E16').Value = '  +4.03%  '

# Row 17
$ws.Range('D17').Value = '2.432.88'
$ws.Range('E17').Value = '  +4.11%  '

# Row 18
$ws.Range('D18').Value = '11.35'
$ws.Range('E18').Value = '  +6.33%  '

# Row 19
$ws.Range('E19').Value = '  +4.02%  '

# Row 20
$ws.Range('D20').Value = '334.08'
$ws.Range('E20').Value = '  +0.51%  '

# Row 21
$ws.Range('E21').Value = '  +5.19%  '

# Row 22
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.18%  '

# Row 23
$ws.Range('D23').Value = '64.65'
$ws.Range('E23').Value = '  +3.02%  '

# Row 24
$ws.Range('E24').Value = '  +0.12%  '

# Row 25
$ws.Range('E25').Value = '  -0.04%  '

# Row 26
$ws.Range('D26').Value = '8.39'
$ws.Range('E26').Value = '  -1.57%  '

# Row 27
$ws.Range('E27').Value = '  -4.31%  '

# Row 28
$ws.Range('E28').Value = '  +1.81%  '

# Row 29
$ws.Range('D29').Value = '0.0₃0769'
$ws.Range('E29').Value = '  +4.95%  '

# Row 30
$ws.Range('D30').Value = '171.27'
$ws.Range('E30').Value = '  +0.41%  '

# Row 31
$ws.Range('D31').Value = '6.24'
$ws.Range('E31').Value = '  +2.48%  '

# Row 32
$ws.Range('D32').Value = '18.70'
$ws.Range('E32').Value = '  +1.10%  '

# Row 33
$ws.Range('D33').Value = '1.02'
$ws.Range('E33').Value = '  +1.13%  '

# Row 34
$ws.Range('E34').Value = '  -0.03%  '

# Row 35
$ws.Range('D35').Value = '4.27'
$ws.Range('E35').Value = '  +1.31%  '

# Row 36
$ws.Range('E36').Value = '  +3.43%  '

# Row 37
$ws.Range('E37').Value = '  +0.12%  '

# Row 38
$ws.Range('D38').Value = '1.63'
$ws.Range('E38').Value = '  -0.43%  '

# Row 39
$ws.Range('D39').Value = '40.22'
$ws.Range('E39').Value = '  +3.35%  '

# Row 40
$ws.Range('D40').Value = '0.422'
$ws.Range('E40').Value = '  +12.62%  '

# Row 41
$ws.Range('D41').Value = '301.75'
$ws.Range('E41').Value = '  +5.68%  '

# Row 42
$ws.Range('D42').Value = '3.73'
$ws.Range('E42').Value = '  +3.00%  '

# Row 43
$ws.Range('D43').Value = '142.25'
$ws.Range('E43').Value = '  -1.57%  '

# Row 44
$ws.Range('D44').Value = '0.0959'
$ws.Range('E44').Value = '  +2.05%  '

# Row 45
$ws.Range('E45').Value = '  +4.25%  '

# Row 46
$ws.Range('D46').Value = '0.412'
$ws.Range('E46').Value = '  +8.01%  '

# Row 47
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = '0.569'
$ws.Range('E47').Value = '  +1.22%  '

# Row 48
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '18.96'
$ws.Range('E48').Value = '  -0.11%  '

# Row 49
$ws.Range('E49').Value = '  +3.57%  '

# Row 50
$ws.Range('D50').Value = '11.05'
$ws.Range('E50').Value = '  -0.22%  '

# Row 51
$ws.Range('D51').Value = '1.57'
$ws.Range('E51').Value = '  +2.77%  '

# Restore the cells to the default "Normal" style now that the text has
# been written, so no stray style/number-format is left behind on them.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
